# repull data, push all data, mean calculation
# Update the dSF (column F) values for the rows that were re-pulled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -6
$ws.Range("F5").Value = -2
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = 3
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = -2
$ws.Range("F13").Value = 8
$ws.Range("F14").Value = -8
$ws.Range("F15").Value = -1
